$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 590, pushing the existing rows (590-649) down to (592-651).
$ws.Rows.Item(590).Resize(2).Insert()

# New row 590: Española / Primera / $/unidad / Llay Llay
$ws.Cells.Item(590,1).Value = 3
$ws.Cells.Item(590,2).Value = "Femacal de La Calera"
$ws.Cells.Item(590,3).Value = "Coquimbo"
$ws.Cells.Item(590,4).Value = 45212
$ws.Cells.Item(590,5).Value = 5
$ws.Cells.Item(590,6).Value = 100112013
$ws.Cells.Item(590,7).Value = "Alcachofa"
$ws.Cells.Item(590,8).Value = "Española"
$ws.Cells.Item(590,9).Value = "Primera"
$ws.Cells.Item(590,10).Value = 22000
$ws.Cells.Item(590,11).Value = 380
$ws.Cells.Item(590,12).Value = 400
$ws.Cells.Item(590,13).Value = 391
$ws.Cells.Item(590,14).Value = "`$/unidad"
$ws.Cells.Item(590,15).Value = "Llay Llay"
$ws.Cells.Item(590,16).Value = 391
$ws.Cells.Item(590,17).Value = 1
$ws.Cells.Item(590,18).Value = "Hortaliza"

# New row 591: Española / Segunda / $/unidad / Llay Llay
$ws.Cells.Item(591,1).Value = 3
$ws.Cells.Item(591,2).Value = "Femacal de La Calera"
$ws.Cells.Item(591,3).Value = "Coquimbo"
$ws.Cells.Item(591,4).Value = 45212
$ws.Cells.Item(591,5).Value = 5
$ws.Cells.Item(591,6).Value = 100112013
$ws.Cells.Item(591,7).Value = "Alcachofa"
$ws.Cells.Item(591,8).Value = "Española"
$ws.Cells.Item(591,9).Value = "Segunda"
$ws.Cells.Item(591,10).Value = 12000
$ws.Cells.Item(591,11).Value = 250
$ws.Cells.Item(591,12).Value = 250
$ws.Cells.Item(591,13).Value = 250
$ws.Cells.Item(591,14).Value = "`$/unidad"
$ws.Cells.Item(591,15).Value = "Llay Llay"
$ws.Cells.Item(591,16).Value = 250
$ws.Cells.Item(591,17).Value = 1
$ws.Cells.Item(591,18).Value = "Hortaliza"

Write-Host "Inserted rows 590-591 with new Alcachofa price data"
